$d = $word.ActiveDocument

# The document has three consecutive placeholder tables (4th, 5th, 6th
# "View daily log within time period" tables). The edit:
#   - deletes the two trailing placeholder tables entirely (along with the
#     blank spacer paragraph that separates each table), and
#   - clears the surviving (4th) table's first-row/first-cell text so only
#     an empty paragraph remains (keeping its second row "Description" /
#     "Acceptance Criteria" content as-is).

# Delete the last ("6th") placeholder table plus the spacer paragraph that
# precedes it.
$t6 = $d.Tables.Item(6)
$t6.Delete()

$t5 = $d.Tables.Item(5)
$spacerStart = $t5.Range.End
$d.Range($spacerStart, $spacerStart + 1).Delete()

# Delete the ("5th") placeholder table plus the spacer paragraph that
# precedes it.
$t5 = $d.Tables.Item(5)
$t5.Delete()

$t4 = $d.Tables.Item(4)
$spacerStart = $t4.Range.End
$d.Range($spacerStart, $spacerStart + 1).Delete()

# Clear the text (and lastRenderedPageBreak run) out of the surviving
# table's first cell, leaving a single empty paragraph behind.
$t4 = $d.Tables.Item(4)
$cell = $t4.Cell(1, 1)
$cellRange = $cell.Range
$d.Range($cellRange.Start, $cellRange.End - 1).Delete()
